# Update cached market-price / profit figures pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 26.5
$ws.Range("J2").Value = 85
$ws.Range("L2").Value = 85
$ws.Range("N2").Value = -311
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H19").Value = 2076.2
$ws.Range("I19").Value = 1439
$ws.Range("K19").Value = 1439
$ws.Range("M19").Value = -1264
$ws.Range("H43").Value = 242436.3
$ws.Range("I43").Value = 1055
$ws.Range("J43").Value = 456997.44
$ws.Range("K43").Value = 1055
$ws.Range("L43").Value = 456997.44
$ws.Range("M43").Value = -986
$ws.Range("N43").Value = -457135.44
$ws.Range("H62").Value = 58835204
$ws.Range("I62").Value = 111112020
$ws.Range("K62").Value = 111112020
$ws.Range("M62").Value = -111111396
$ws.Range("H64").Value = 7747.25
$ws.Range("J64").Value = 7994.5
$ws.Range("L64").Value = 7994.5
$ws.Range("N64").Value = -8490.5
$ws.Range("H65").Value = 58835204
$ws.Range("I65").Value = 111112020
$ws.Range("K65").Value = 555560100
$ws.Range("M65").Value = -555556980
$ws.Range("H67").Value = 7747.25
$ws.Range("J67").Value = 7994.5
$ws.Range("L67").Value = 7994.5
$ws.Range("N67").Value = -9710.5
$ws.Range("H87").Value = 31109.777
$ws.Range("J87").Value = 31109.777
$ws.Range("L87").Value = 31109.777
$ws.Range("N87").Value = -33605.777
$ws.Range("H90").Value = 31109.777
$ws.Range("J90").Value = 31109.777
$ws.Range("L90").Value = 93329.33099999999
$ws.Range("N90").Value = -105809.331
$ws.Range("H137").Value = 2656.2917
$ws.Range("I137").Value = 2511.9524
$ws.Range("J137").Value = 3666.6667
$ws.Range("K137").Value = 7535.8572
$ws.Range("L137").Value = 11000.0001
$ws.Range("M137").Value = -4985.8572
$ws.Range("N137").Value = -16100.0001
$ws.Range("H138").Value = 3993.2222
$ws.Range("I138").Value = 1138.9487
$ws.Range("J138").Value = 7366.4546
$ws.Range("K138").Value = 3416.8461
$ws.Range("L138").Value = 22099.3638
$ws.Range("M138").Value = 1723.1539
$ws.Range("N138").Value = -32379.3638
$ws.Range("H141").Value = 7247965
$ws.Range("I141").Value = 7577236
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 22731708
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -22726528
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2900.5264
$ws.Range("I2").Value = 988.75
$ws.Range("K2").Value = 988.75
$ws.Range("M2").Value = -875.75
$ws.Range("H32").Value = 1670918.1
$ws.Range("I32").Value = 1764484
$ws.Range("K32").Value = 1764484
$ws.Range("M32").Value = -1764197
$ws.Range("H45").Value = 11098.25
$ws.Range("I45").Value = 3474.5
$ws.Range("K45").Value = 3474.5
$ws.Range("M45").Value = -3097.5
$ws.Range("H61").Value = 4698.1333
$ws.Range("I61").Value = 2551.875
$ws.Range("K61").Value = 2551.875
$ws.Range("M61").Value = -2339.875
$ws.Range("H74").Value = 18224.135
$ws.Range("I74").Value = 23454.445
$ws.Range("K74").Value = 23454.445
$ws.Range("M74").Value = -22580.445
$ws.Range("H77").Value = 18224.135
$ws.Range("I77").Value = 23454.445
$ws.Range("K77").Value = 117272.225
$ws.Range("M77").Value = -112904.225
$ws.Range("H116").Value = 2900.5264
$ws.Range("I116").Value = 988.75
$ws.Range("K116").Value = 988.75
$ws.Range("M116").Value = 1305.25
$ws.Range("H122").Value = 28644.334
$ws.Range("I122").Value = 68333.336
$ws.Range("J122").Value = 8799.833000000001
$ws.Range("K122").Value = 205000.008
$ws.Range("L122").Value = 26399.499
$ws.Range("M122").Value = -202550.008
$ws.Range("N122").Value = -31299.499
$ws.Range("H136").Value = 4698.1333
$ws.Range("I136").Value = 2551.875
$ws.Range("K136").Value = 7655.625
$ws.Range("M136").Value = -5105.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2900.5264
$ws.Range("I3").Value = 988.75
$ws.Range("K3").Value = 988.75
$ws.Range("M3").Value = -874.75
$ws.Range("H26").Value = 29992
$ws.Range("I26").Value = 29992
$ws.Range("K26").Value = 29992
$ws.Range("M26").Value = -29700
$ws.Range("H96").Value = 6500
$ws.Range("I96").Value = 6500
$ws.Range("K96").Value = 6500
$ws.Range("M96").Value = -3754
$ws.Range("H107").Value = 40181964
$ws.Range("I107").Value = 59213400
$ws.Range("J107").Value = 4490.3335
$ws.Range("K107").Value = 59213400
$ws.Range("L107").Value = 4490.3335
$ws.Range("M107").Value = -59211480
$ws.Range("N107").Value = -8330.333500000001
$ws.Range("H134").Value = 4703.418
$ws.Range("I134").Value = 1326.7667
$ws.Range("K134").Value = 3980.300099999999
$ws.Range("M134").Value = -1445.300099999999
$ws.Range("H138").Value = 84955.60000000001
$ws.Range("J138").Value = 84955.60000000001
$ws.Range("L138").Value = 84955.60000000001
$ws.Range("N138").Value = -95235.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8200900.5
$ws.Range("I58").Value = 12196373
$ws.Range("J58").Value = 10181.6
$ws.Range("K58").Value = 12196373
$ws.Range("L58").Value = 10181.6
$ws.Range("M58").Value = -12196170
$ws.Range("N58").Value = -10587.6
$ws.Range("H105").Value = 2647923.2
$ws.Range("I105").Value = 3106236
$ws.Range("K105").Value = 3106236
$ws.Range("M105").Value = -3104489
$ws.Range("H132").Value = 4681.78
$ws.Range("I132").Value = 2019.625
$ws.Range("K132").Value = 6058.875
$ws.Range("M132").Value = -3528.875
$ws.Range("H134").Value = 6123.6284
$ws.Range("I134").Value = 1407.2941
$ws.Range("K134").Value = 4221.8823
$ws.Range("M134").Value = -1686.8823
$ws.Range("H136").Value = 8200900.5
$ws.Range("I136").Value = 12196373
$ws.Range("J136").Value = 10181.6
$ws.Range("K136").Value = 36589119
$ws.Range("L136").Value = 30544.8
$ws.Range("M136").Value = -36586569
$ws.Range("N136").Value = -35644.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3614.4443
$ws.Range("J68").Value = 5264.2
$ws.Range("L68").Value = 15792.6
$ws.Range("N68").Value = -17414.6
$ws.Range("H71").Value = 3614.4443
$ws.Range("J71").Value = 5264.2
$ws.Range("L71").Value = 47377.8
$ws.Range("N71").Value = -55489.8
$ws.Range("H107").Value = 14286075
$ws.Range("J107").Value = 18182168
$ws.Range("L107").Value = 54546504
$ws.Range("N107").Value = -54550344

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2783.44
$ws.Range("I102").Value = 2592.7144
$ws.Range("J102").Value = 3784.75
$ws.Range("K102").Value = 2592.7144
$ws.Range("L102").Value = 3784.75
$ws.Range("M102").Value = -970.7143999999998
$ws.Range("N102").Value = -7028.75
$ws.Range("H122").Value = 5971030
$ws.Range("J122").Value = 3749.25
$ws.Range("L122").Value = 11247.75
$ws.Range("N122").Value = -16147.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8875.625
$ws.Range("I40").Value = 8000
$ws.Range("K40").Value = 8000
$ws.Range("M40").Value = -7864
$ws.Range("H61").Value = 5016.619
$ws.Range("I61").Value = 1595.3334
$ws.Range("J61").Value = 6385.1333
$ws.Range("K61").Value = 1595.3334
$ws.Range("L61").Value = 6385.1333
$ws.Range("M61").Value = -1393.3334
$ws.Range("N61").Value = -6789.1333
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H113").Value = 5016.619
$ws.Range("I113").Value = 1595.3334
$ws.Range("J113").Value = 6385.1333
$ws.Range("K113").Value = 1595.3334
$ws.Range("L113").Value = 6385.1333
$ws.Range("M113").Value = 574.6666
$ws.Range("N113").Value = -10725.1333
$ws.Range("H122").Value = 7586.6665
$ws.Range("I122").Value = 6465.8335
$ws.Range("J122").Value = 8333.888999999999
$ws.Range("K122").Value = 19397.5005
$ws.Range("L122").Value = 25001.667
$ws.Range("M122").Value = -16947.5005
$ws.Range("N122").Value = -29901.667
$ws.Range("H136").Value = 11330.346
$ws.Range("I136").Value = 4498.75
$ws.Range("J136").Value = 12572.454
$ws.Range("K136").Value = 13496.25
$ws.Range("L136").Value = 37717.362
$ws.Range("M136").Value = -10946.25
$ws.Range("N136").Value = -42817.362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 200580.19
$ws.Range("J122").Value = 6891.5835
$ws.Range("L122").Value = 20674.7505
$ws.Range("N122").Value = -25574.7505
$ws.Range("H136").Value = 25671956
$ws.Range("I136").Value = 62500772
$ws.Range("K136").Value = 187502316
$ws.Range("M136").Value = -187499766

